# Auto-generated edit script applying the commit diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4503
$ws.Range("F4").Value = 439
$ws.Range("F5").Value = 3608
$ws.Range("F6").Value = 1037
$ws.Range("F7").Value = 164
$ws.Range("F9").Value = 351
$ws.Range("F10").Value = 350
$ws.Range("F11").Value = 2487
$ws.Range("F12").Value = 1274
$ws.Range("F14").Value = 1969
$ws.Range("F20").Value = 10293
$ws.Range("F21").Value = 6007
$ws.Range("F22").Value = 11
$ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202408/bGrekRm71722496057963.jpeg"
$ws.Range("C23").Value = "杭州·首届次元之门动漫游戏博览会懒喵N²次元·爱内里菜日间专场签售"
$ws.Range("F23").Value = 5
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202408/rCJwjRU51722495672321.jpeg"
$ws.Range("F24").Value = 391
$ws.Range("F25").Value = 212
$ws.Range("C26").Value = "杭州·首届次元之门动漫游戏博览会懒喵N²次元·NANOナノ日间专场签售"
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202408/wFo0KVEW1722495461964.jpeg"
$ws.Range("C27").Value = "杭州·首届次元之门动漫游戏博览会懒喵N²次元·Survive Said The Prophet日间专场签售"
$ws.Range("F27").Value = 8
$ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202408/zU4ao5PO1722495828005.jpeg"
$ws.Range("F28").Value = 834
$ws.Range("F30").Value = 165
$ws.Range("F31").Value = 852
$ws.Range("F32").Value = 3555
$ws.Range("F35").Value = 472
$ws.Range("F37").Value = 255
$ws.Range("F39").Value = 237
$ws.Range("G40").Value = 39
$ws.Range("F42").Value = 1115
$ws.Range("F43").Value = 159
$ws.Range("F44").Value = 58
$ws.Range("F45").Value = 86
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1
$ws.Range("F7").Value = 14
$ws.Range("F12").Value = 132
$ws.Range("F15").Value = 3553
$ws.Range("F16").Value = 79
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8775
$ws.Range("F4").Value = 1610
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1610
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 4503
$ws.Range("F7").Value = 439
$ws.Range("F8").Value = 3609
$ws.Range("F9").Value = 1037
$ws.Range("F10").Value = 164
$ws.Range("F12").Value = 350
$ws.Range("F13").Value = 2487
$ws.Range("F15").Value = 1274
$ws.Range("F20").Value = 132
$ws.Range("F24").Value = 10293
$ws.Range("F25").Value = 3553
$ws.Range("F26").Value = 79
$ws.Range("F27").Value = 11
$ws.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202408/bGrekRm71722496057963.jpeg"
$ws.Range("F28").Value = 391
$ws.Range("F29").Value = 212
$ws.Range("C30").Value = "杭州·首届次元之门动漫游戏博览会懒喵N²次元·NANOナノ日间专场签售"
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202408/wFo0KVEW1722495461964.jpeg"
$ws.Range("C31").Value = "杭州·首届次元之门动漫游戏博览会懒喵N²次元·Survive Said The Prophet日间专场签售"
$ws.Range("F31").Value = 8
$ws.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202408/zU4ao5PO1722495828005.jpeg"
$ws.Range("F32").Value = 834
$ws.Range("F34").Value = 165
$ws.Range("F35").Value = 852
$ws.Range("F36").Value = 3555
$ws.Range("F39").Value = 255
$ws.Range("F41").Value = 237
$ws.Range("G42").Value = 39
$ws.Range("F44").Value = 1115
$ws.Range("F45").Value = 159
$ws.Range("F46").Value = 86
